$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.672.10"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.631.63"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.23"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "1.859.60"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "1.648.01"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "26.656.22"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.45"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +8.13%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.82"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.85"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "1.217.92"
$ws.Range("E36").Value = "  +5.05%  "
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D44").Value = "1.769.39"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.76"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.06"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.64"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.52%  "
